$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 14-17 (Resolving-Mac is no longer a receiving/target cluster;
# the underlying TPM recompute dropped these rows entirely)
$ws.Range("A14:T17").Delete()

# Update remaining data rows (2-13) with recomputed TPM-based values
$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Wnt4'
$ws.Range("C2").Value = 'Fzd6'
$ws.Range("D2").Value = 'ECs'
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.646881333333333
$ws.Range("H2").Value = 4.940644
$ws.Range("I2").Value = 0.3367300927127475
$ws.Range("J2").Value = 0.3367300927127475
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 28.58650933333334
$ws.Range("N2").Value = 85.759528
$ws.Range("O2").Value = 0.9059490896276022
$ws.Range("P2").Value = 0.9059490896276023
$ws.Range("Q2").Value = 47.07858860622578
$ws.Range("R2").Value = 423.707297456032
$ws.Range("S2").Value = 0.3050603209433317
$ws.Range("T2").Value = 0.3050603209433317
$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Wnt4'
$ws.Range("C3").Value = 'Fzd6'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.646881333333333
$ws.Range("H3").Value = 4.940644
$ws.Range("I3").Value = 0.3367300927127475
$ws.Range("J3").Value = 0.3367300927127475
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.662736333333334
$ws.Range("N3").Value = 7.988209000000001
$ws.Range("O3").Value = 0.08438608327351124
$ws.Range("P3").Value = 0.08438608327351126
$ws.Range("Q3").Value = 4.385210762955111
$ws.Range("R3").Value = 39.466896866596
$ws.Range("S3").Value = 0.02841533364435507
$ws.Range("T3").Value = 0.02841533364435508
$ws.Range("A4").Value = 'ECs'
$ws.Range("B4").Value = 'Wnt4'
$ws.Range("C4").Value = 'Fzd6'
$ws.Range("D4").Value = 'MuSCs'
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.646881333333333
$ws.Range("H4").Value = 4.940644
$ws.Range("I4").Value = 0.3367300927127475
$ws.Range("J4").Value = 0.3367300927127475
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.304966
$ws.Range("N4").Value = 0.914898
$ws.Range("O4").Value = 0.009664827098886481
$ws.Range("P4").Value = 0.009664827098886483
$ws.Range("Q4").Value = 0.5022428127013333
$ws.Range("R4").Value = 4.520185314312
$ws.Range("S4").Value = 0.003254438125060719
$ws.Range("T4").Value = 0.00325443812506072
$ws.Range("A5").Value = 'FAPs'
$ws.Range("B5").Value = 'Wnt4'
$ws.Range("C5").Value = 'Fzd6'
$ws.Range("D5").Value = 'ECs'
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.533623
$ws.Range("H5").Value = 4.600869
$ws.Range("I5").Value = 0.3135726931406526
$ws.Range("J5").Value = 0.3135726931406525
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 28.58650933333334
$ws.Range("N5").Value = 85.759528
$ws.Range("O5").Value = 0.9059490896276022
$ws.Range("P5").Value = 0.9059490896276023
$ws.Range("Q5").Value = 43.84092820331468
$ws.Range("R5").Value = 394.5683538298321
$ws.Range("S5").Value = 0.2840808958828497
$ws.Range("T5").Value = 0.2840808958828496
$ws.Range("A6").Value = 'FAPs'
$ws.Range("B6").Value = 'Wnt4'
$ws.Range("C6").Value = 'Fzd6'
$ws.Range("D6").Value = 'FAPs'
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.533623
$ws.Range("H6").Value = 4.600869
$ws.Range("I6").Value = 0.3135726931406526
$ws.Range("J6").Value = 0.3135726931406525
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.662736333333334
$ws.Range("N6").Value = 7.988209000000001
$ws.Range("O6").Value = 0.08438608327351124
$ws.Range("P6").Value = 0.08438608327351126
$ws.Range("Q6").Value = 4.083633683735668
$ws.Range("R6").Value = 36.75270315362101
$ws.Range("S6").Value = 0.02646117139566629
$ws.Range("T6").Value = 0.02646117139566629
$ws.Range("A7").Value = 'FAPs'
$ws.Range("B7").Value = 'Wnt4'
$ws.Range("C7").Value = 'Fzd6'
$ws.Range("D7").Value = 'MuSCs'
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.533623
$ws.Range("H7").Value = 4.600869
$ws.Range("I7").Value = 0.3135726931406526
$ws.Range("J7").Value = 0.3135726931406525
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.304966
$ws.Range("N7").Value = 0.914898
$ws.Range("O7").Value = 0.009664827098886481
$ws.Range("P7").Value = 0.009664827098886483
$ws.Range("Q7").Value = 0.4677028718180001
$ws.Range("R7").Value = 4.209325846362
$ws.Range("S7").Value = 0.003030625862136594
$ws.Range("T7").Value = 0.003030625862136594
$ws.Range("A8").Value = 'MuSCs'
$ws.Range("B8").Value = 'Wnt4'
$ws.Range("C8").Value = 'Fzd6'
$ws.Range("D8").Value = 'ECs'
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.115861333333333
$ws.Range("H8").Value = 3.347584
$ws.Range("I8").Value = 0.2281549269050169
$ws.Range("J8").Value = 0.2281549269050168
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 28.58650933333334
$ws.Range("N8").Value = 85.759528
$ws.Range("O8").Value = 0.9059490896276022
$ws.Range("P8").Value = 0.9059490896276023
$ws.Range("Q8").Value = 31.89858042003911
$ws.Range("R8").Value = 287.087223780352
$ws.Range("S8").Value = 0.2066967483236522
$ws.Range("T8").Value = 0.2066967483236522
$ws.Range("A9").Value = 'MuSCs'
$ws.Range("B9").Value = 'Wnt4'
$ws.Range("C9").Value = 'Fzd6'
$ws.Range("D9").Value = 'FAPs'
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.115861333333333
$ws.Range("H9").Value = 3.347584
$ws.Range("I9").Value = 0.2281549269050169
$ws.Range("J9").Value = 0.2281549269050168
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.662736333333334
$ws.Range("N9").Value = 7.988209000000001
$ws.Range("O9").Value = 0.08438608327351124
$ws.Range("P9").Value = 0.08438608327351126
$ws.Range("Q9").Value = 2.971244515228445
$ws.Range("R9").Value = 26.741200637056
$ws.Range("S9").Value = 0.01925310066106863
$ws.Range("T9").Value = 0.01925310066106863
$ws.Range("A10").Value = 'MuSCs'
$ws.Range("B10").Value = 'Wnt4'
$ws.Range("C10").Value = 'Fzd6'
$ws.Range("D10").Value = 'MuSCs'
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.115861333333333
$ws.Range("H10").Value = 3.347584
$ws.Range("I10").Value = 0.2281549269050169
$ws.Range("J10").Value = 0.2281549269050168
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.304966
$ws.Range("N10").Value = 0.914898
$ws.Range("O10").Value = 0.009664827098886481
$ws.Range("P10").Value = 0.009664827098886483
$ws.Range("Q10").Value = 0.3402997673813334
$ws.Range("R10").Value = 3.062697906432
$ws.Range("S10").Value = 0.002205077920296072
$ws.Range("T10").Value = 0.002205077920296072
$ws.Range("A11").Value = 'Resolving-Mac'
$ws.Range("B11").Value = 'Wnt4'
$ws.Range("C11").Value = 'Fzd6'
$ws.Range("D11").Value = 'ECs'
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5944396666666667
$ws.Range("H11").Value = 1.783319
$ws.Range("I11").Value = 0.1215422872415831
$ws.Range("J11").Value = 0.1215422872415831
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 28.58650933333334
$ws.Range("N11").Value = 85.759528
$ws.Range("O11").Value = 0.9059490896276022
$ws.Range("P11").Value = 0.9059490896276023
$ws.Range("Q11").Value = 16.99295507927022
$ws.Range("R11").Value = 152.936595713432
$ws.Range("S11").Value = 0.1101111244777688
$ws.Range("T11").Value = 0.1101111244777688
$ws.Range("A12").Value = 'Resolving-Mac'
$ws.Range("B12").Value = 'Wnt4'
$ws.Range("C12").Value = 'Fzd6'
$ws.Range("D12").Value = 'FAPs'
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5944396666666667
$ws.Range("H12").Value = 1.783319
$ws.Range("I12").Value = 0.1215422872415831
$ws.Range("J12").Value = 0.1215422872415831
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.662736333333334
$ws.Range("N12").Value = 7.988209000000001
$ws.Range("O12").Value = 0.08438608327351124
$ws.Range("P12").Value = 0.08438608327351126
$ws.Range("Q12").Value = 1.582836098407889
$ws.Range("R12").Value = 14.245524885671
$ws.Range("S12").Value = 0.01025647757242126
$ws.Range("T12").Value = 0.01025647757242126
$ws.Range("A13").Value = 'Resolving-Mac'
$ws.Range("B13").Value = 'Wnt4'
$ws.Range("C13").Value = 'Fzd6'
$ws.Range("D13").Value = 'MuSCs'
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5944396666666667
$ws.Range("H13").Value = 1.783319
$ws.Range("I13").Value = 0.1215422872415831
$ws.Range("J13").Value = 0.1215422872415831
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.304966
$ws.Range("N13").Value = 0.914898
$ws.Range("O13").Value = 0.009664827098886481
$ws.Range("P13").Value = 0.009664827098886483
$ws.Range("Q13").Value = 0.1812838873846667
$ws.Range("R13").Value = 1.631554986462
$ws.Range("S13").Value = 0.001174685191393097
$ws.Range("T13").Value = 0.001174685191393098
